$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.976.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.037.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.032.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("E13").Value = "  +2.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.81%  "

$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.537.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.929.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.034.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.92%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("E34").Value = "  +2.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0875"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.08%  "

$ws.Range("E36").Value = "  +1.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "

$ws.Range("E38").Value = "  +10.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.309"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "399.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.733.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.43%  "
